$d = $word.ActiveDocument

$d.Content.Find.Execute("3+57=60", $true, $false, $false, $false, $false, $true, 1, $false, "88+9=97", 2) | Out-Null
$d.Content.Find.Execute("5+75=80", $true, $false, $false, $false, $false, $true, 1, $false, "57-24=33", 2) | Out-Null
$d.Content.Find.Execute("13+78=91", $true, $false, $false, $false, $false, $true, 1, $false, "97-83=14", 2) | Out-Null
$d.Content.Find.Execute("90-82=8", $true, $false, $false, $false, $false, $true, 1, $false, "10+19=29", 2) | Out-Null
$d.Content.Find.Execute("68+30=98", $true, $false, $false, $false, $false, $true, 1, $false, "11+30=41", 2) | Out-Null
$d.Content.Find.Execute("33-17=16", $true, $false, $false, $false, $false, $true, 1, $false, "22-21=1", 2) | Out-Null
$d.Content.Find.Execute("18+30=48", $true, $false, $false, $false, $false, $true, 1, $false, "59-1=58", 2) | Out-Null
$d.Content.Find.Execute("47-33=14", $true, $false, $false, $false, $false, $true, 1, $false, "83-49=34", 2) | Out-Null
$d.Content.Find.Execute("40+9=49", $true, $false, $false, $false, $false, $true, 1, $false, "1+0=1", 2) | Out-Null
$d.Content.Find.Execute("3+48=51", $true, $false, $false, $false, $false, $true, 1, $false, "2+87=89", 2) | Out-Null
$d.Content.Find.Execute("23-7=16", $true, $false, $false, $false, $false, $true, 1, $false, "61+0=61", 2) | Out-Null
$d.Content.Find.Execute("15-2=13", $true, $false, $false, $false, $false, $true, 1, $false, "88-75=13", 2) | Out-Null
$d.Content.Find.Execute("36+16=52", $true, $false, $false, $false, $false, $true, 1, $false, "77-27=50", 2) | Out-Null
$d.Content.Find.Execute("39+29=68", $true, $false, $false, $false, $false, $true, 1, $false, "79+13=92", 2) | Out-Null
$d.Content.Find.Execute("84+1=85", $true, $false, $false, $false, $false, $true, 1, $false, "14+56=70", 2) | Out-Null
$d.Content.Find.Execute("71-49=22", $true, $false, $false, $false, $false, $true, 1, $false, "58+39=97", 2) | Out-Null
$d.Content.Find.Execute("94-76=18", $true, $false, $false, $false, $false, $true, 1, $false, "23-21=2", 2) | Out-Null
$d.Content.Find.Execute("81-68=13", $true, $false, $false, $false, $false, $true, 1, $false, "13+24=37", 2) | Out-Null
$d.Content.Find.Execute("39+15=54", $true, $false, $false, $false, $false, $true, 1, $false, "97-72=25", 2) | Out-Null
$d.Content.Find.Execute("49+5=54", $true, $false, $false, $false, $false, $true, 1, $false, "67-3=64", 2) | Out-Null
$d.Content.Find.Execute("11+5=16", $true, $false, $false, $false, $false, $true, 1, $false, "68+12=80", 2) | Out-Null
$d.Content.Find.Execute("10+1=11", $true, $false, $false, $false, $false, $true, 1, $false, "50+27=77", 2) | Out-Null
$d.Content.Find.Execute("82-25=57", $true, $false, $false, $false, $false, $true, 1, $false, "32+59=91", 2) | Out-Null
$d.Content.Find.Execute("11+4=15", $true, $false, $false, $false, $false, $true, 1, $false, "15-0=15", 2) | Out-Null
$d.Content.Find.Execute("84-56=28", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=68", 2) | Out-Null
$d.Content.Find.Execute("67+7=74", $true, $false, $false, $false, $false, $true, 1, $false, "25+65=90", 2) | Out-Null
$d.Content.Find.Execute("82-33=49", $true, $false, $false, $false, $false, $true, 1, $false, "68-67=1", 2) | Out-Null
$d.Content.Find.Execute("69-48=21", $true, $false, $false, $false, $false, $true, 1, $false, "70+25=95", 2) | Out-Null
$d.Content.Find.Execute("80-70=10", $true, $false, $false, $false, $false, $true, 1, $false, "79-11=68", 2) | Out-Null
$d.Content.Find.Execute("24+65=89", $true, $false, $false, $false, $false, $true, 1, $false, "74-73=1", 2) | Out-Null
$d.Content.Find.Execute("10+88=98", $true, $false, $false, $false, $false, $true, 1, $false, "27-9=18", 2) | Out-Null
$d.Content.Find.Execute("86-9=77", $true, $false, $false, $false, $false, $true, 1, $false, "44+28=72", 2) | Out-Null
$d.Content.Find.Execute("94-62=32", $true, $false, $false, $false, $false, $true, 1, $false, "86-26=60", 2) | Out-Null
$d.Content.Find.Execute("55-39=16", $true, $false, $false, $false, $false, $true, 1, $false, "26+69=95", 2) | Out-Null
$d.Content.Find.Execute("38+52=90", $true, $false, $false, $false, $false, $true, 1, $false, "3+56=59", 2) | Out-Null
$d.Content.Find.Execute("85-39=46", $true, $false, $false, $false, $false, $true, 1, $false, "49-35=14", 2) | Out-Null
$d.Content.Find.Execute("47-6=41", $true, $false, $false, $false, $false, $true, 1, $false, "47-31=16", 2) | Out-Null
$d.Content.Find.Execute("23-19=4", $true, $false, $false, $false, $false, $true, 1, $false, "78+11=89", 2) | Out-Null
$d.Content.Find.Execute("43-33=10", $true, $false, $false, $false, $false, $true, 1, $false, "80+14=94", 2) | Out-Null
$d.Content.Find.Execute("89-28=61", $true, $false, $false, $false, $false, $true, 1, $false, "36+20=56", 2) | Out-Null
$d.Content.Find.Execute("8+77=85", $true, $false, $false, $false, $false, $true, 1, $false, "76+11=87", 2) | Out-Null
$d.Content.Find.Execute("50-7=43", $true, $false, $false, $false, $false, $true, 1, $false, "73-20=53", 2) | Out-Null
$d.Content.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=47", 2) | Out-Null
$d.Content.Find.Execute("49+6=55", $true, $false, $false, $false, $false, $true, 1, $false, "99-25=74", 2) | Out-Null
$d.Content.Find.Execute("71-63=8", $true, $false, $false, $false, $false, $true, 1, $false, "33-28=5", 2) | Out-Null
$d.Content.Find.Execute("60-45=15", $true, $false, $false, $false, $false, $true, 1, $false, "64-60=4", 2) | Out-Null
$d.Content.Find.Execute("86+13=99", $true, $false, $false, $false, $false, $true, 1, $false, "87-28=59", 2) | Out-Null
$d.Content.Find.Execute("53+22=75", $true, $false, $false, $false, $false, $true, 1, $false, "20+49=69", 2) | Out-Null
$d.Content.Find.Execute("22+9=31", $true, $false, $false, $false, $false, $true, 1, $false, "11+55=66", 2) | Out-Null
$d.Content.Find.Execute("7+12=19", $true, $false, $false, $false, $false, $true, 1, $false, "80+9=89", 2) | Out-Null
$d.Content.Find.Execute("36-21=15", $true, $false, $false, $false, $false, $true, 1, $false, "24-15=9", 2) | Out-Null
$d.Content.Find.Execute("47-17=30", $true, $false, $false, $false, $false, $true, 1, $false, "0+87=87", 2) | Out-Null
$d.Content.Find.Execute("12+9=21", $true, $false, $false, $false, $false, $true, 1, $false, "62-55=7", 2) | Out-Null
$d.Content.Find.Execute("63+31=94", $true, $false, $false, $false, $false, $true, 1, $false, "15+71=86", 2) | Out-Null
$d.Content.Find.Execute("34+60=94", $true, $false, $false, $false, $false, $true, 1, $false, "66+18=84", 2) | Out-Null
$d.Content.Find.Execute("57+28=85", $true, $false, $false, $false, $false, $true, 1, $false, "28-11=17", 2) | Out-Null
$d.Content.Find.Execute("13+64=77", $true, $false, $false, $false, $false, $true, 1, $false, "28+19=47", 2) | Out-Null
$d.Content.Find.Execute("73-51=22", $true, $false, $false, $false, $false, $true, 1, $false, "97-14=83", 2) | Out-Null
$d.Content.Find.Execute("2+67=69", $true, $false, $false, $false, $false, $true, 1, $false, "92-72=20", 2) | Out-Null
$d.Content.Find.Execute("70-29=41", $true, $false, $false, $false, $false, $true, 1, $false, "87-17=70", 2) | Out-Null
$d.Content.Find.Execute("54+5=59", $true, $false, $false, $false, $false, $true, 1, $false, "6+58=64", 2) | Out-Null
$d.Content.Find.Execute("45-1=44", $true, $false, $false, $false, $false, $true, 1, $false, "65+29=94", 2) | Out-Null
$d.Content.Find.Execute("8+23=31", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=39", 2) | Out-Null
$d.Content.Find.Execute("46-33=13", $true, $false, $false, $false, $false, $true, 1, $false, "87-59=28", 2) | Out-Null
$d.Content.Find.Execute("51+6=57", $true, $false, $false, $false, $false, $true, 1, $false, "95-45=50", 2) | Out-Null
$d.Content.Find.Execute("6-3=3", $true, $false, $false, $false, $false, $true, 1, $false, "29+24=53", 2) | Out-Null
$d.Content.Find.Execute("92-26=66", $true, $false, $false, $false, $false, $true, 1, $false, "63+4=67", 2) | Out-Null
$d.Content.Find.Execute("22+38=60", $true, $false, $false, $false, $false, $true, 1, $false, "69-59=10", 2) | Out-Null
$d.Content.Find.Execute("35+32=67", $true, $false, $false, $false, $false, $true, 1, $false, "72-43=29", 2) | Out-Null
$d.Content.Find.Execute("27+58=85", $true, $false, $false, $false, $false, $true, 1, $false, "29-26=3", 2) | Out-Null
$d.Content.Find.Execute("28+36=64", $true, $false, $false, $false, $false, $true, 1, $false, "79+16=95", 2) | Out-Null
$d.Content.Find.Execute("14+70=84", $true, $false, $false, $false, $false, $true, 1, $false, "57+21=78", 2) | Out-Null
$d.Content.Find.Execute("3+73=76", $true, $false, $false, $false, $false, $true, 1, $false, "43+20=63", 2) | Out-Null
$d.Content.Find.Execute("34-20=14", $true, $false, $false, $false, $false, $true, 1, $false, "96-94=2", 2) | Out-Null
$d.Content.Find.Execute("4-2=2", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=74", 2) | Out-Null
$d.Content.Find.Execute("81-13=68", $true, $false, $false, $false, $false, $true, 1, $false, "32+32=64", 2) | Out-Null
$d.Content.Find.Execute("90+6=96", $true, $false, $false, $false, $false, $true, 1, $false, "69+14=83", 2) | Out-Null
$d.Content.Find.Execute("0+13=13", $true, $false, $false, $false, $false, $true, 1, $false, "84-46=38", 2) | Out-Null
$d.Content.Find.Execute("23+73=96", $true, $false, $false, $false, $false, $true, 1, $false, "50-40=10", 2) | Out-Null
$d.Content.Find.Execute("12+44=56", $true, $false, $false, $false, $false, $true, 1, $false, "53-28=25", 2) | Out-Null
$d.Content.Find.Execute("1+35=36", $true, $false, $false, $false, $false, $true, 1, $false, "47+3=50", 2) | Out-Null
$d.Content.Find.Execute("97-66=31", $true, $false, $false, $false, $false, $true, 1, $false, "77-20=57", 2) | Out-Null
$d.Content.Find.Execute("19+55=74", $true, $false, $false, $false, $false, $true, 1, $false, "54+6=60", 2) | Out-Null
$d.Content.Find.Execute("86+6=92", $true, $false, $false, $false, $false, $true, 1, $false, "54-36=18", 2) | Out-Null
$d.Content.Find.Execute("38-20=18", $true, $false, $false, $false, $false, $true, 1, $false, "44+6=50", 2) | Out-Null
$d.Content.Find.Execute("16-14=2", $true, $false, $false, $false, $false, $true, 1, $false, "47-47=0", 2) | Out-Null
$d.Content.Find.Execute("19-16=3", $true, $false, $false, $false, $false, $true, 1, $false, "52+37=89", 2) | Out-Null
$d.Content.Find.Execute("2+22=24", $true, $false, $false, $false, $false, $true, 1, $false, "35+52=87", 2) | Out-Null
$d.Content.Find.Execute("64-57=7", $true, $false, $false, $false, $false, $true, 1, $false, "14+37=51", 2) | Out-Null
$d.Content.Find.Execute("1+64=65", $true, $false, $false, $false, $false, $true, 1, $false, "76-3=73", 2) | Out-Null
$d.Content.Find.Execute("49-26=23", $true, $false, $false, $false, $false, $true, 1, $false, "37+42=79", 2) | Out-Null
$d.Content.Find.Execute("26-16=10", $true, $false, $false, $false, $false, $true, 1, $false, "84-75=9", 2) | Out-Null
$d.Content.Find.Execute("42+38=80", $true, $false, $false, $false, $false, $true, 1, $false, "18+79=97", 2) | Out-Null
$d.Content.Find.Execute("7+32=39", $true, $false, $false, $false, $false, $true, 1, $false, "88-40=48", 2) | Out-Null
$d.Content.Find.Execute("13+77=90", $true, $false, $false, $false, $false, $true, 1, $false, "21+26=47", 2) | Out-Null
$d.Content.Find.Execute("77+0=77", $true, $false, $false, $false, $false, $true, 1, $false, "8+67=75", 2) | Out-Null
$d.Content.Find.Execute("85-21=64", $true, $false, $false, $false, $false, $true, 1, $false, "62-11=51", 2) | Out-Null
$d.Content.Find.Execute("34+18=52", $true, $false, $false, $false, $false, $true, 1, $false, "37+25=62", 2) | Out-Null
$d.Content.Find.Execute("76-30=46", $true, $false, $false, $false, $false, $true, 1, $false, "12+0=12", 2) | Out-Null
$d.Content.Find.Execute("82-40=42", $true, $false, $false, $false, $false, $true, 1, $false, "85-61=24", 2) | Out-Null
